# "run prepare & render with final data"
#
# The simulation was re-run with a new country (Russia) added and all
# survey-share figures recomputed. This inserts a "Russia" data column
# between Japan and Saudi Arabia (shifting Saudi Arabia -> M, USA -> N),
# refreshes every numeric share in the table with the newly computed
# values, rewraps two of the long label strings, and swaps the
# "Could sign a petition..." row label for the new
# "Would support a global movement..." label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:N1) : countries -----------------------------------
$ws.Range("B1").Value2 = '$ bold(''All'')'
$ws.Range("C1").Value2 = '$ bold(''Europe'')'
$ws.Range("D1").Value2 = 'France'
$ws.Range("E1").Value2 = 'Germany'
$ws.Range("F1").Value2 = 'Italy'
$ws.Range("G1").Value2 = 'Poland'
$ws.Range("H1").Value2 = 'Spain'
$ws.Range("I1").Value2 = 'United Kingdom'
$ws.Range("J1").Value2 = 'Switzerland'
$ws.Range("K1").Value2 = 'Japan'
$ws.Range("L1").Value2 = 'Russia'
$ws.Range("M1").Value2 = 'Saudi Arabia'
$ws.Range("N1").Value2 = 'USA'

# --- Row labels (A2:A10) : survey questions ----------------------------
$ws.Range("A2").Value2 = 'Supports tax on world top 1% to finance global poverty reduction
(Additional 15% tax on income over [$120k/year in PPP])'
$ws.Range("A3").Value2 = 'Supports tax on world top 3% to finance global poverty reduction
(Additional 15% tax over [$80k], 30% over [$120k], 45% over [$1M])'
$ws.Range("A4").Value2 = 'Prefers sustainable future*'
$ws.Range("A5").Value2 = '"Governments should actively cooperate to have all countries
converge in terms of GDP per capita by the end of the century"'
$ws.Range("A6").Value2 = 'Would support a global movement to tackle CC, tax millionaires,
 and fund LICs (either petition, demonstrate, strike, or donate)'
$ws.Range("A7").Value2 = 'More likely to vote for party if part of worldwide
coalition for climate action and global redistribution'
$ws.Range("A8").Value2 = 'Supports reparations for colonization and slavery in
the form of funding education and technology transfers'
$ws.Range("A9").Value2 = '"My taxes should go towards solving global problems"'
$ws.Range("A10").Value2 = '"My taxes ... global problems" (Global Nation, 2024)'

# --- Data cells : refreshed shares for every row/country --------------

# Row 2 - Supports tax on world top 1%...
$ws.Range("B2").Value2 = 0.693996120467564
$ws.Range("C2").Value2 = 0.727267133682238
$ws.Range("D2").Value2 = 0.713266166753258
$ws.Range("E2").Value2 = 0.722752212925052
$ws.Range("F2").Value2 = 0.839866307760825
$ws.Range("G2").Value2 = 0.694516301386687
$ws.Range("H2").Value2 = 0.727791884744065
$ws.Range("I2").Value2 = 0.671225049697161
$ws.Range("J2").Value2 = 0.603572620182626
$ws.Range("K2").Value2 = 0.691681993781983
$ws.Range("L2").Value2 = 0.754855071318372
$ws.Range("M2").Value2 = 0.817345461694808
$ws.Range("N2").Value2 = 0.622335501686223

# Row 3 - Supports tax on world top 3%...
$ws.Range("B3").Value2 = 0.641188431291775
$ws.Range("C3").Value2 = 0.663594295554775
$ws.Range("D3").Value2 = 0.698772999587592
$ws.Range("E3").Value2 = 0.622101654744799
$ws.Range("F3").Value2 = 0.706881256898736
$ws.Range("G3").Value2 = 0.699228862448248
$ws.Range("H3").Value2 = 0.658500586296013
$ws.Range("I3").Value2 = 0.665804066558109
$ws.Range("J3").Value2 = 0.420104611305175
$ws.Range("K3").Value2 = 0.550394111458431
$ws.Range("L3").Value2 = 0.756622372880961
$ws.Range("M3").Value2 = 0.816920760799125
$ws.Range("N3").Value2 = 0.57421952327369

# Row 4 - Prefers sustainable future*
$ws.Range("B4").Value2 = 0.680881448179833
$ws.Range("C4").Value2 = 0.701539116816613
$ws.Range("D4").Value2 = 0.719673932778947
$ws.Range("E4").Value2 = 0.704645953003574
$ws.Range("F4").Value2 = 0.760434225154626
$ws.Range("G4").Value2 = 0.572739428737827
$ws.Range("H4").Value2 = 0.73608837047601
$ws.Range("I4").Value2 = 0.679200342355269
$ws.Range("J4").Value2 = 0.665468662233116
$ws.Range("K4").Value2 = 0.758076861129753
$ws.Range("L4").Value2 = 0.688615273248795
$ws.Range("M4").Value2 = 0.713280127381035
$ws.Range("N4").Value2 = 0.616918649447641

# Row 5 - "Governments should actively cooperate..."
$ws.Range("B5").Value2 = 0.704965329416964
$ws.Range("C5").Value2 = 0.780310385878786
$ws.Range("D5").Value2 = 0.761687470249255
$ws.Range("E5").Value2 = 0.758489172765931
$ws.Range("F5").Value2 = 0.874100026122781
$ws.Range("G5").Value2 = 0.843596425922759
$ws.Range("H5").Value2 = 0.842976707711552
$ws.Range("I5").Value2 = 0.658792898536381
$ws.Range("J5").Value2 = 0.660378064524565
$ws.Range("K5").Value2 = 0.703727559204962
$ws.Range("L5").Value2 = 0.778963825426238
$ws.Range("M5").Value2 = 0.930231790695484
$ws.Range("N5").Value2 = 0.561007628053552

# Row 6 - Would support a global movement...
$ws.Range("B6").Value2 = 0.675595447215337
$ws.Range("C6").Value2 = 0.719216740354837
$ws.Range("D6").Value2 = 0.699222514786681
$ws.Range("E6").Value2 = 0.688082663981164
$ws.Range("F6").Value2 = 0.819338712934373
$ws.Range("G6").Value2 = 0.708420268414952
$ws.Range("H6").Value2 = 0.74352822863702
$ws.Range("I6").Value2 = 0.68138828161491
$ws.Range("J6").Value2 = 0.639531813440066
$ws.Range("K6").Value2 = 0.557841849059486
$ws.Range("L6").ClearContents()
$ws.Range("M6").Value2 = 0.727098526374066
$ws.Range("N6").Value2 = 0.666651932459956

# Row 7 - More likely to vote for party...
$ws.Range("B7").Value2 = 0.682631646934764
$ws.Range("C7").Value2 = 0.724046548532723
$ws.Range("D7").Value2 = 0.71331522240695
$ws.Range("E7").Value2 = 0.70786414424213
$ws.Range("F7").Value2 = 0.817943755181141
$ws.Range("G7").Value2 = 0.642319300405405
$ws.Range("H7").Value2 = 0.765355431118643
$ws.Range("I7").Value2 = 0.695445383221243
$ws.Range("J7").Value2 = 0.578158388278922
$ws.Range("K7").Value2 = 0.557598270145635
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value2 = 0.669431535311024

# Row 8 - Supports reparations for colonization...
$ws.Range("B8").Value2 = 0.451810364536854
$ws.Range("C8").Value2 = 0.502618407568093
$ws.Range("D8").Value2 = 0.437375697519216
$ws.Range("E8").Value2 = 0.43866370099415
$ws.Range("F8").Value2 = 0.696737544675098
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value2 = 0.507401270303301
$ws.Range("I8").Value2 = 0.458595587653966
$ws.Range("J8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value2 = 0.401158464849325

# Row 9 - "My taxes should go towards solving global problems"
$ws.Range("B9").Value2 = 0.592122368373113
$ws.Range("C9").Value2 = 0.609896497471146
$ws.Range("D9").Value2 = 0.432040689538255
$ws.Range("E9").Value2 = 0.620111679474376
$ws.Range("F9").Value2 = 0.766975357539963
$ws.Range("G9").Value2 = 0.618309495622478
$ws.Range("H9").Value2 = 0.708786744917504
$ws.Range("I9").Value2 = 0.574354517212933
$ws.Range("J9").Value2 = 0.538340937683521
$ws.Range("K9").Value2 = 0.584137117430032
$ws.Range("L9").Value2 = 0.573182308826315
$ws.Range("M9").Value2 = 0.888820570273345
$ws.Range("N9").Value2 = 0.554219568199052

# Row 10 - "My taxes ... global problems" (Global Nation, 2024) -- all #NUM!
$ws.Range("B10").Value2 = '#NUM!'
$ws.Range("C10").Value2 = '#NUM!'
$ws.Range("D10").Value2 = '#NUM!'
$ws.Range("E10").Value2 = '#NUM!'
$ws.Range("F10").Value2 = '#NUM!'
$ws.Range("G10").Value2 = '#NUM!'
$ws.Range("H10").Value2 = '#NUM!'
$ws.Range("I10").Value2 = '#NUM!'
$ws.Range("J10").ClearContents()
$ws.Range("K10").Value2 = '#NUM!'
$ws.Range("L10").ClearContents()
$ws.Range("M10").Value2 = '#NUM!'
$ws.Range("N10").Value2 = '#NUM!'
